$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 243 (old rows 243-247 shift down to 244-248).
$ws.Rows.Item(243).Insert()

# Populate the newly inserted row 243 with the new weekly data point.
$ws.Cells.Item(243, 1).Value = 7
$ws.Cells.Item(243, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(243, 3).Value = "Ñuble"
$ws.Cells.Item(243, 4).Value = 44448
$ws.Cells.Item(243, 5).Value = 16
$ws.Cells.Item(243, 6).Value = 100114014
$ws.Cells.Item(243, 7).Value = "Betarraga"
$ws.Cells.Item(243, 8).Value = "Sin especificar"
$ws.Cells.Item(243, 9).Value = "Primera"
$ws.Cells.Item(243, 10).Value = 300
$ws.Cells.Item(243, 11).Value = 750
$ws.Cells.Item(243, 12).Value = 800
$ws.Cells.Item(243, 13).Value = 775
$ws.Cells.Item(243, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(243, 15).Value = "Región del Maule"
$ws.Cells.Item(243, 16).Value = 155
$ws.Cells.Item(243, 17).Value = 5
$ws.Cells.Item(243, 18).Value = "Hortaliza"
